# Add team record (Wins / Losses / Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new header cells (AD1:AF1) the same formatting as the existing
# header row (bold font + border), copying it from the last header cell (AC1)
# before writing the new header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

# Header row (row 1): new column headers in AD, AE, AF.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-53: same team record repeated for every player.
$lastRow = 53
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 68   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 94   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
